# Update Name of Algo
# Applies updated RandomForest imputation values to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 8.947999999999992
$ws.Range("D4").Value = -7.641200000000001
$ws.Range("B7").Value = 5.511200000000003
$ws.Range("D12").Value = -5.813599999999999
$ws.Range("B16").Value = 5.070199999999998
$ws.Range("D18").Value = -8.992899999999993
$ws.Range("D19").Value = -8.544599999999994
$ws.Range("D20").Value = -8.542399999999988
$ws.Range("B28").Value = 5.958200000000004
$ws.Range("B29").Value = 5.166800000000004
$ws.Range("D31").Value = -7.414899999999994
$ws.Range("B32").Value = 6.502599999999997
$ws.Range("B40").Value = 9.179900000000002
$ws.Range("D40").Value = -8.865599999999995
$ws.Range("D42").Value = -8.667899999999998
$ws.Range("D47").Value = -7.588499999999997
$ws.Range("D48").Value = -7.403199999999997
$ws.Range("B52").Value = 5.245899999999997
$ws.Range("B57").Value = 4.790299999999996
$ws.Range("D63").Value = -6.728899999999997
$ws.Range("D64").Value = -7.288199999999993
$ws.Range("B66").Value = 5.809699999999996
$ws.Range("D76").Value = -7.8152
$ws.Range("D81").Value = -7.6144
$ws.Range("D89").Value = -8.290299999999997
$ws.Range("D94").Value = -5.942899999999998
$ws.Range("B100").Value = 4.858699999999999

$wb.Save()
